# Add "user 5" results row to the userstudy results sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 ("user 4"): fill in the missing "M.C./?" marker in column B
$ws.Cells.Item(6, 2).Value = "?"

# Row 7 ("user 5"): newly added results
$ws.Cells.Item(7, 2).Value = "?"    # B7 - Moving C marker (case 1)
$ws.Cells.Item(7, 3).Value = 39     # C7 - case 1 (bp) steps
$ws.Cells.Item(7, 4).Value = 318    # D7 - case 1 (qp) time(s)
$ws.Cells.Item(7, 9).Value = 3      # I7 - case 2 (qp) steps
$ws.Cells.Item(7, 10).Value = 195   # J7 - case 2 (qp) time(s)

# Move the active selection to reflect the newly entered data
$ws.Range("B8").Select()
